$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh adds two new price records (dated 2023-02-10, serial 44967)
# at the top of the data block, pushing the existing rows (15-26) down to
# rows 17-28. Mirror that with a native row insert so formatting/styles of
# the surrounding rows (e.g. the date format on column D) carry over.
$ws.Rows("15:16").Insert()

# Row 15: Terminal Hortofrutícola Agro Chillán - Arándano (blue) - Primera
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C15").Value = "Ñuble"
$ws.Range("D15").Value = 44967
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100101
$ws.Range("H15").Value = "Berries"
$ws.Range("I15").Value = 100101001
$ws.Range("J15").Value = "Arándano (blue)"
$ws.Range("K15").Value = "Sin especificar"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 3000
$ws.Range("O15").Value = 3000
$ws.Range("P15").Value = 3000
$ws.Range("Q15").Value = "$/bandeja 2 kilos"
$ws.Range("R15").Value = "Provincia de Diguillín"
$ws.Range("S15").Value = 1500
$ws.Range("T15").Value = 2

# Row 16: Terminal Hortofrutícola Agro Chillán - Arándano (blue) - Segunda
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 44967
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100101
$ws.Range("H16").Value = "Berries"
$ws.Range("I16").Value = 100101001
$ws.Range("J16").Value = "Arándano (blue)"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 2500
$ws.Range("O16").Value = 2500
$ws.Range("P16").Value = 2500
$ws.Range("Q16").Value = "$/bandeja 2 kilos"
$ws.Range("R16").Value = "Provincia de Diguillín"
$ws.Range("S16").Value = 1250
$ws.Range("T16").Value = 2
